# Update Name of Algo
# Applies updated values produced by a re-run of the RandomForest imputation
# algorithm for a handful of cells on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value  = 6.144799999999988
$ws.Range("E5").Value  = 12.68379999999999
$ws.Range("E9").Value  = 13.30160000000001
$ws.Range("E11").Value = 13.89449999999999
$ws.Range("B21").Value = 5.780199999999999
$ws.Range("E21").Value = 12.78069999999999
$ws.Range("B23").Value = 5.614300000000002
$ws.Range("B25").Value = 6.062499999999991
